$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The post "「アラビア語が好き」" (row 299) was removed from the blog export.
# Delete that entire row; Excel shifts every row below it up by one,
# which matches all of the following rows renumbering down by 1
# (300->299, 301->300, ... 387->386) with their content unchanged.
$ws.Rows.Item(299).Delete()
